$d = $word.ActiveDocument

# Helper characters that must not go through Word's "smart quotes" /
# AutoFormat-as-you-type autocorrection (Find.Execute's Replace argument
# gets auto-corrected, so every substitution below locates the match with
# Find.Execute and then assigns Range.Text directly, which is not
# autocorrected).
$rsquo = [char]0x2019
$euro  = [char]0x20AC

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText)
    if ($found) {
        $rng.Text = $newText
    }
    return $found
}

# 1. Title / heading (appears twice: the Heading1 and the bold line near
#    the bottom of the document) - replace both occurrences.
$oldTitle = "Play Big Size Fishin' for Free - Slot Game Review"
$newTitle = "Play Big Size Fishin" + $rsquo + " Free: Review of Features and Gameplay"
Replace-Text $oldTitle $newTitle
Replace-Text $oldTitle $newTitle

# 2. "What we like" bullet list
Replace-Text "Easy to understand gameplay" "Visually stunning with crisp lines and color shading"
Replace-Text "Attractive RTP for experienced gamblers" "Immersive music"
Replace-Text "Visually stunning with crisp lines and colors" "Convenient Autoplay feature"
Replace-Text "Free demo version available" "Medium-high variance with attractive RTP"

# 3. "What we don't like" bullet list
Replace-Text "Only 10 fixed paylines" "Limited number of paylines"
Replace-Text ("Limited maximum bet of " + $euro + "50") "Only available at online casinos"

# 4. Meta description (italic text at the very end of the document)
$oldMeta = "Explore the fishing adventure of Big Size Fishin" + $rsquo + " slot game. Review of gameplay, symbols, and free spins. Try the demo version before playing."
$newMeta = "Discover the gameplay and features of Big Size Fishin' and play it for free. Review of Big Size Fishin' slot game."
Replace-Text $oldMeta $newMeta
